$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the sheet from "Sheet1" to "Tiles" ---
$ws.Name = "Tiles"

# --- 2. Grow the sheet from 44 to 52 data rows. ---
#     Insert the 8 new rows by copying the last existing row (44) and using
#     Insert so each new row inherits the correct cell styles (s="1"/s="2")
#     instead of ending up with no style, the way a bare Cells.Item write would.
for ($r = 45; $r -le 52; $r++) {
    $ws.Rows.Item($r - 1).Copy()
    $ws.Rows.Item($r).Insert(-4121)   # -4121 = xlShiftDown
}
$excel.CutCopyMode = 0

# --- 3. Write the new tile-map contents for columns A, B and C. ---
#     Column A: row index (0 on row1, then 1,1,1,2,2,2,3,3,3,4,4,4,5,5,6, then +1 each row)
#     Column B: sub-index cycling 1/2/3 for the first rows, then always 1
#     Column C: constant 2
$colA = @(0,1,1,2,2,2,3,3,3,4,4,4,5,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43)
$colB = @(1,1,2,1,2,3,1,2,3,1,2,3,1,2,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1)
$colC = @(2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2)

for ($r = 1; $r -le 52; $r++) {
    $idx = $r - 1

    if ($r -ge 16) {
        $ws.Cells.Item($r, 1).Formula = "=A" + ($r - 1) + "+1"
    } else {
        $ws.Cells.Item($r, 1).Value = $colA[$idx]
    }

    $ws.Cells.Item($r, 2).Value = $colB[$idx]
    $ws.Cells.Item($r, 3).Value = $colC[$idx]
}

# --- 4. Rows 7-11 pick up a slightly smaller explicit row height in the saved file. ---
for ($r = 7; $r -le 11; $r++) {
    $ws.Rows.Item($r).RowHeight = 13.8
}

# --- 5. Recreate the selection shown in the saved file: C2:C52 active at C2. ---
[void]$ws.Range("C2:C52").Select()
